$wb = $excel.ActiveWorkbook

# --- Matches_SOG: append new match rows 447-451 ---
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

# row 447
$wsMatches.Range("A447").NumberFormat = "@"
$wsMatches.Range("A447").Value = "897744"
$wsMatches.Range("A447").Style = "Normal"
$wsMatches.Range("B447").Value = "2025-11-08T17:00:00"
$wsMatches.Range("C447").Value = "Драконы"
$wsMatches.Range("D447").Value = "Северсталь"
$wsMatches.Range("E447").Value = 39
$wsMatches.Range("F447").Value = 36
$wsMatches.Range("G447").Value = "khl_text"

# row 448
$wsMatches.Range("A448").NumberFormat = "@"
$wsMatches.Range("A448").Value = "897745"
$wsMatches.Range("A448").Style = "Normal"
$wsMatches.Range("B448").Value = "2025-11-08T17:00:00"
$wsMatches.Range("C448").Value = "Локомотив"
$wsMatches.Range("D448").Value = "Ак Барс"
$wsMatches.Range("E448").Value = 27
$wsMatches.Range("F448").Value = 29
$wsMatches.Range("G448").Value = "khl_text"

# row 449
$wsMatches.Range("A449").NumberFormat = "@"
$wsMatches.Range("A449").Value = "897746"
$wsMatches.Range("A449").Style = "Normal"
$wsMatches.Range("B449").Value = "2025-11-08T17:00:00"
$wsMatches.Range("C449").Value = "ЦСКА"
$wsMatches.Range("D449").Value = "Сибирь"
$wsMatches.Range("E449").Value = 31
$wsMatches.Range("F449").Value = 24
$wsMatches.Range("G449").Value = "khl_text"

# row 450
$wsMatches.Range("A450").NumberFormat = "@"
$wsMatches.Range("A450").Value = "897748"
$wsMatches.Range("A450").Style = "Normal"
$wsMatches.Range("B450").Value = "2025-11-08T19:30:00"
$wsMatches.Range("C450").Value = "Динамо М"
$wsMatches.Range("D450").Value = "СКА"
$wsMatches.Range("E450").Value = 30
$wsMatches.Range("F450").Value = 31
$wsMatches.Range("G450").Value = "khl_text"

# row 451
$wsMatches.Range("A451").NumberFormat = "@"
$wsMatches.Range("A451").Value = "897747"
$wsMatches.Range("A451").Style = "Normal"
$wsMatches.Range("B451").Value = "2025-11-08T17:10:00"
$wsMatches.Range("C451").Value = "Динамо Мн"
$wsMatches.Range("D451").Value = "Лада"
$wsMatches.Range("E451").Value = 49
$wsMatches.Range("F451").Value = 14
$wsMatches.Range("G451").Value = "khl_text"

# --- Shots_HA updates ---
$wsHA = $wb.Worksheets.Item("Shots_HA")
$wsHA.Range("D2").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D3").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D4").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D5").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("F5").Value = 20
$wsHA.Range("K5").Value = 659
$wsHA.Range("L5").Value = 595
$wsHA.Range("M5").Value = 33
$wsHA.Range("N5").Value = 29.8
$wsHA.Range("D6").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D7").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D8").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("E8").Value = 17
$wsHA.Range("G8").Value = 563
$wsHA.Range("H8").Value = 470
$wsHA.Range("I8").Value = 33.1
$wsHA.Range("J8").Value = 27.6
$wsHA.Range("D9").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("E9").Value = 23
$wsHA.Range("G9").Value = 836
$wsHA.Range("H9").Value = 617
$wsHA.Range("I9").Value = 36.3
$wsHA.Range("J9").Value = 26.8
$wsHA.Range("D10").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("E10").Value = 20
$wsHA.Range("G10").Value = 571
$wsHA.Range("H10").Value = 698
$wsHA.Range("I10").Value = 28.6
$wsHA.Range("J10").Value = 34.9
$wsHA.Range("D11").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("F11").Value = 19
$wsHA.Range("K11").Value = 495
$wsHA.Range("L11").Value = 723
$wsHA.Range("M11").Value = 26.1
$wsHA.Range("N11").Value = 38.1
$wsHA.Range("D12").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("E12").Value = 18
$wsHA.Range("G12").Value = 555
$wsHA.Range("H12").Value = 495
$wsHA.Range("I12").Value = 30.8
$wsHA.Range("J12").Value = 27.5
$wsHA.Range("D13").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D14").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D15").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("F15").Value = 16
$wsHA.Range("K15").Value = 499
$wsHA.Range("L15").Value = 517
$wsHA.Range("N15").Value = 32.3
$wsHA.Range("D16").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D17").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("F17").Value = 24
$wsHA.Range("K17").Value = 795
$wsHA.Range("L17").Value = 638
$wsHA.Range("M17").Value = 33.1
$wsHA.Range("N17").Value = 26.6
$wsHA.Range("D18").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("F18").Value = 22
$wsHA.Range("K18").Value = 612
$wsHA.Range("L18").Value = 681
$wsHA.Range("M18").Value = 27.8
$wsHA.Range("D19").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D20").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D21").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D22").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("D23").Value = "2025-11-08T19:30:00Z"
$wsHA.Range("E23").Value = 19
$wsHA.Range("G23").Value = 443
$wsHA.Range("H23").Value = 543
$wsHA.Range("I23").Value = 23.3
$wsHA.Range("J23").Value = 28.6

# --- Shots_Summary updates ---
$wsSummary = $wb.Worksheets.Item("Shots_Summary")
$wsSummary.Range("D2").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D3").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D4").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D5").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("E5").Value = 43
$wsSummary.Range("F5").Value = 1438
$wsSummary.Range("G5").Value = 1182
$wsSummary.Range("H5").Value = 33.4
$wsSummary.Range("D6").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D7").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D8").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("E8").Value = 38
$wsSummary.Range("F8").Value = 1159
$wsSummary.Range("G8").Value = 1143
$wsSummary.Range("D9").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("E9").Value = 41
$wsSummary.Range("F9").Value = 1504
$wsSummary.Range("G9").Value = 1101
$wsSummary.Range("H9").Value = 36.7
$wsSummary.Range("I9").Value = 26.9
$wsSummary.Range("D10").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("E10").Value = 40
$wsSummary.Range("F10").Value = 1123
$wsSummary.Range("G10").Value = 1435
$wsSummary.Range("H10").Value = 28.1
$wsSummary.Range("D11").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("E11").Value = 41
$wsSummary.Range("F11").Value = 1093
$wsSummary.Range("G11").Value = 1499
$wsSummary.Range("H11").Value = 26.7
$wsSummary.Range("I11").Value = 36.6
$wsSummary.Range("D12").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("E12").Value = 43
$wsSummary.Range("F12").Value = 1332
$wsSummary.Range("G12").Value = 1103
$wsSummary.Range("H12").Value = 31
$wsSummary.Range("I12").Value = 25.7
$wsSummary.Range("D13").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D14").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D15").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("E15").Value = 41
$wsSummary.Range("F15").Value = 1329
$wsSummary.Range("G15").Value = 1351
$wsSummary.Range("H15").Value = 32.4
$wsSummary.Range("D16").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D17").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("E17").Value = 40
$wsSummary.Range("F17").Value = 1264
$wsSummary.Range("G17").Value = 1001
$wsSummary.Range("H17").Value = 31.6
$wsSummary.Range("I17").Value = 25
$wsSummary.Range("D18").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("E18").Value = 41
$wsSummary.Range("F18").Value = 1123
$wsSummary.Range("G18").Value = 1414
$wsSummary.Range("H18").Value = 27.4
$wsSummary.Range("I18").Value = 34.5
$wsSummary.Range("D19").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D20").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D21").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D22").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("D23").Value = "2025-11-08T19:30:00Z"
$wsSummary.Range("E23").Value = 40
$wsSummary.Range("F23").Value = 966
$wsSummary.Range("G23").Value = 1140
$wsSummary.Range("H23").Value = 24.1
$wsSummary.Range("I23").Value = 28.5

# --- Meta_ext updates ---
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = "2025-11-08T19:30:00Z"
$wsMeta.Range("D2").Value = 54
